$wb = $excel.ActiveWorkbook

$wsServer = $wb.Worksheets.Item("ServerDatabaseConfig")
$wsScripts = $wb.Worksheets.Item("ScriptsConfig")

# Duplicate the last data row (row 4) down into rows 5 and 6 so the new
# rows pick up the same cell formatting (fills/styles) as the existing
# "SQLScript" rows, then overwrite the values that actually change.
$wsScripts.Range("A4:G4").Copy($wsScripts.Range("A5:G5")) | Out-Null
$wsScripts.Range("A4:G4").Copy($wsScripts.Range("A6:G6")) | Out-Null

# Row 5: RecentQuries.sql
$wsScripts.Range("A5").Value = 1
$wsScripts.Range("B5").Value = "SQLScript"
$wsScripts.Range("C5").Value = "C:\Z_Tests\SQLScripts"
$wsScripts.Range("D5").Value = "RecentQuries.sql"
$wsScripts.Range("E5").Value = '$var1|var1_replace;$var2|var2_replace'
$wsScripts.Range("F5").Value = 1
$wsScripts.Range("G5").Value = 15

# Row 6: RecentQuriesByUsers.sql
$wsScripts.Range("A6").Value = 1
$wsScripts.Range("B6").Value = "SQLScript"
$wsScripts.Range("C6").Value = "C:\Z_Tests\SQLScripts"
$wsScripts.Range("D6").Value = "RecentQuriesByUsers.sql"
$wsScripts.Range("E6").Value = '$var1|var1_replace;$var2|var2_replace'
$wsScripts.Range("F6").Value = 1
$wsScripts.Range("G6").Value = 15

# Update the selection left behind on ServerDatabaseConfig ...
$wsServer.Range("D10").Select() | Out-Null

# ... then switch focus to ScriptsConfig, which becomes the active tab.
$wsScripts.Activate() | Out-Null
$wsScripts.Range("C11").Select() | Out-Null
